$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bugs")

# Sorting finished, icons added:
#  - "Absent sorting" bug (row 4) is now In Review
#  - "Add icons" bug (row 6) is now In Review
#  - "Add info panel" bug (row 10) is now In Progress
$ws.Range("C4").Value = "Review"
$ws.Range("C6").Value = "Review"
$ws.Range("C10").Value = "InProgress"

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("F21").Select()
